$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full refreshed leaderboard data (rows 2-22), columns A-J:
# Rang, Speler, Score, 180'ers, 100+ finishes, Totaal Score, Aantal Darts, 3-Darts Gemiddelde, Totaal, Winnaar
$data = @(
    @(1, "Yannick den Daggelder", 42, 1, 0, 27236, 1542, 52.99, 43, 2),
    @(2, "Remco Riem", 37, 4, 0, 15478, 713, 65.12, 41, 2),
    @(3, "Alessandro Delia", 26, 0, 0, 17627, 1033, 51.19, 26, 0),
    @(3, "Nick Fitzpatrick", 24, 1, 1, 11558, 544, 63.74, 26, 1),
    @(5, "Dartin Dan", 17, 1, 1, 12667, 722, 52.63, 19, 0),
    @(6, "Louis Tweddle", 17, 0, 1, 11864, 581, 61.26, 18, 0),
    @(6, "Rocky Van Den Eeckhoudt", 16, 2, 0, 11771, 632, 55.88, 18, 0),
    @(8, "Gijs Tromp", 14, 1, 0, 11451, 624, 55.05, 15, 0),
    @(9, "Robin Willis", 14, 0, 0, 9658, 601, 48.21, 14, 0),
    @(10, "Milan Schoenmakers", 11, 1, 0, 6984, 333, 62.92, 12, 0),
    @(11, "Max Walter", 11, 0, 0, 6262, 477, 39.38, 11, 0),
    @(11, "Afendi Kelana", 8, 1, 2, 8344, 432, 57.94, 11, 0),
    @(13, "Diego Meerveld", 9, 0, 0, 7022, 457, 46.1, 9, 0),
    @(13, "Francesco Petruccelli", 8, 1, 0, 6404, 349, 55.05, 9, 0),
    @(13, "martin Berry", 8, 1, 0, 4295, 285, 45.21, 9, 0),
    @(16, "Magnus Gladh", 7, 0, 0, 6401, 377, 50.94, 7, 0),
    @(17, "Tristan Snoep", 6, 0, 0, 6258, 462, 40.64, 6, 0),
    @(18, "Oscar Ebbeling", 3, 0, 0, 3730, 258, 43.37, 3, 0),
    @(19, "Niels van Dommelen", 2, 0, 0, 2266, 135, 50.36, 2, 0),
    @(19, "Noah B", 2, 0, 0, 2410, 140, 51.64, 2, 0),
    @(21, "Lee Dolan", 1, 0, 0, 1770, 131, 40.53, 1, 0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
}
